# PAGOPROVEEDORES.xlsx edit
# "agregar correo para Bancolom Davivien y CA, CC para bancolombia"
#
# Rows 110-119 of Hoja1 get their Banco / Tipo de Producto / Referencia filled in
# (columns E, F, G) and the "Tipo de Identificacion" cell in column A is highlighted
# in yellow to flag the rows that were updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# r -> (CodigoBanco, TipoProducto, Referencia)
$updates = @(
    @{ Row = 110; Banco = 51; Tipo = "CC"; Ref = "51000000" },
    @{ Row = 111; Banco = 7;  Tipo = "CC"; Ref = "51000000" },
    @{ Row = 112; Banco = 7;  Tipo = "CA"; Ref = "51000000" },
    @{ Row = 113; Banco = 7;  Tipo = "CC"; Ref = "51000000" },
    @{ Row = 114; Banco = 1;  Tipo = "CC"; Ref = "51000000" },
    @{ Row = 115; Banco = 7;  Tipo = "CC"; Ref = "51000000" },
    @{ Row = 116; Banco = 51; Tipo = "CA"; Ref = "51000000" },
    @{ Row = 117; Banco = 7;  Tipo = "CC"; Ref = "51000000" },
    @{ Row = 118; Banco = 7;  Tipo = "CC"; Ref = "51000000" },
    @{ Row = 119; Banco = 7;  Tipo = "CA"; Ref = "51000000" }
)

foreach ($u in $updates) {
    $r = $u.Row

    $ws.Range("E$r").Value = $u.Banco
    $ws.Range("F$r").Value = $u.Tipo
    $ws.Range("G$r").Value = $u.Ref

    # Highlight the "Tipo de Identificacion" cell (column A) in yellow for the
    # rows that just got completed.
    $ws.Range("A$r").Interior.Color = 65535
}

# Restore the selection to where the user left off working (row 113, column E).
$ws.Range("E113").Select()
